# Auto-generated cell updates for cryptos.xlsx refresh
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue {
    param($CellRef, $Val)
    $r = $ws.Range($CellRef)
    $r.NumberFormat = "@"
    $r.Value = $Val
    $r.Style = "Normal"
}

Set-TextValue "D2" "23.551.80"
Set-TextValue "E2" "  +1.33%  "
Set-TextValue "D3" "1.656.92"
Set-TextValue "E3" "  +2.71%  "
Set-TextValue "E4" "  -0.04%  "
Set-TextValue "D5" "1.000"
Set-TextValue "E5" "  -0.13%  "
Set-TextValue "D6" "302.17"
Set-TextValue "E6" "  -0.12%  "
Set-TextValue "D7" "0.3835"
Set-TextValue "E7" "  +1.28%  "
Set-TextValue "D8" "51.14"
Set-TextValue "E8" "  -1.81%  "
Set-TextValue "D9" "0.3595"
Set-TextValue "E9" "  +1.55%  "
Set-TextValue "D10" "1.240"
Set-TextValue "E10" "  +2.70%  "
Set-TextValue "D11" "0.08190"
Set-TextValue "E11" "  +0.99%  "
Set-TextValue "E12" "  -0.05%  "
Set-TextValue "D13" "22.41"
Set-TextValue "E13" "  +1.29%  "
Set-TextValue "D14" "6.497"
Set-TextValue "E14" "  +1.97%  "
Set-TextValue "D15" "7.502"
Set-TextValue "E15" "  +3.36%  "
Set-TextValue "D16" "0.00001223"
Set-TextValue "E16" "  +1.08%  "
Set-TextValue "D17" "1.650.94"
Set-TextValue "E17" "  +2.08%  "
Set-TextValue "D18" "97.44"
Set-TextValue "E18" "  +3.27%  "
Set-TextValue "D19" "0.06980"
Set-TextValue "E19" "  +0.91%  "
Set-TextValue "D20" "6.829"
Set-TextValue "E20" "  +4.89%  "
Set-TextValue "D21" "17.69"
Set-TextValue "E21" "  +2.77%  "
Set-TextValue "D22" "1.000"
Set-TextValue "E22" "  -0.10%  "
Set-TextValue "D23" "12.71"
Set-TextValue "E23" "  +2.94%  "
Set-TextValue "D24" "23.569.20"
Set-TextValue "E24" "  +1.43%  "
Set-TextValue "D25" "2.499"
Set-TextValue "E25" "  -0.42%  "
Set-TextValue "D26" "2.995"
Set-TextValue "E26" "  -1.20%  "
Set-TextValue "D27" "21.22"
Set-TextValue "E27" "  +1.53%  "
Set-TextValue "D28" "151.97"
Set-TextValue "E28" "  +0.58%  "
Set-TextValue "E29" "  +0.17%  "
Set-TextValue "D30" "133.86"
Set-TextValue "E30" "  +1.01%  "
Set-TextValue "B31" "Filecoin"
Set-TextValue "C31" "https://coinranking.com/coin/ymQub4fuB+filecoin-fil"
Set-TextValue "D31" "7.241"
Set-TextValue "E31" "  +11.63%  "
Set-TextValue "B32" "WrappedliquidstakedEther2.0"
Set-TextValue "C32" "https://coinranking.com/coin/CiixT63n3+wrappedliquidstakedether20-wsteth"
Set-TextValue "D32" "1.840.77"
Set-TextValue "E32" "  +2.88%  "
Set-TextValue "D33" "2.245"
Set-TextValue "E33" "  +7.11%  "
Set-TextValue "D34" "12.04"
Set-TextValue "E34" "  +5.97%  "
Set-TextValue "E35" "  -1.71%  "
Set-TextValue "D36" "0.02801"
Set-TextValue "E36" "  +3.33%  "
Set-TextValue "D37" "6.132"
Set-TextValue "E37" "  +4.85%  "
Set-TextValue "B38" "Stellar"
Set-TextValue "C38" "https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm"
Set-TextValue "D38" "0.08797"
Set-TextValue "E38" "  +0.55%  "
Set-TextValue "B39" "Algorand"
Set-TextValue "C39" "https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo"
Set-TextValue "D39" "0.2497"
Set-TextValue "E39" "  +1.70%  "
Set-TextValue "D40" "0.07027"
Set-TextValue "E40" "  +1.37%  "
Set-TextValue "D41" "13.16"
Set-TextValue "E41" "  +9.86%  "
Set-TextValue "D42" "0.7029"
Set-TextValue "E42" "  +2.06%  "
Set-TextValue "D43" "1.335"
Set-TextValue "E43" "  +0.68%  "
Set-TextValue "D44" "16.10"
Set-TextValue "E44" "  +5.02%  "
Set-TextValue "D45" "0.6547"
Set-TextValue "E45" "  +3.69%  "
Set-TextValue "E46" "  -0.05%  "
Set-TextValue "D47" "2.310"
Set-TextValue "E47" "  +2.74%  "
Set-TextValue "D48" "3.960"
Set-TextValue "E48" "  +0.39%  "
Set-TextValue "D49" "0.07930"
Set-TextValue "E49" "  +0.86%  "
Set-TextValue "D50" "127.88"
Set-TextValue "E50" "  +0.43%  "
Set-TextValue "D51" "1.193"
Set-TextValue "E51" "  +1.84%  "
